$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '25.765.40'
$ws.Range('E2').Value2 = '  -0.20%  '
$ws.Range('D3').Value2 = '1.634.07'
$ws.Range('E3').Value2 = '  -0.12%  '
$ws.Range('E4').Value2 = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '215.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value2 = '  +0.10%  '
$ws.Range('E6').Value2 = '  -0.77%  '
$ws.Range('E7').Value2 = '  -0.15%  '
$ws.Range('E8').Value2 = '  -0.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.0634'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value2 = '  -1.48%  '
$ws.Range('E10').Value2 = '  -1.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.0791'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value2 = '  +0.80%  '
$ws.Range('E12').Value2 = '  +0.27%  '
$ws.Range('D13').Value2 = '1.860.51'
$ws.Range('D14').Value2 = '1.641.90'
$ws.Range('E14').Value2 = '  -0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.560'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value2 = '  +0.44%  '
$ws.Range('E16').Value2 = '  -1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '63.15'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value2 = '  +0.02%  '
$ws.Range('D18').Value2 = '25.781.50'
$ws.Range('E18').Value2 = '  -0.20%  '
$ws.Range('E20').Value2 = '  +0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '192.39'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value2 = '  -0.96%  '
$ws.Range('E22').Value2 = '  +0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '6.29'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value2 = '  +2.28%  '
$ws.Range('E24').Value2 = '  -0.13%  '
$ws.Range('E25').Value2 = '  +2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '142.36'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value2 = '  +2.17%  '
$ws.Range('E27').Value2 = '  +1.81%  '
$ws.Range('E28').Value2 = '  +0.63%  '
$ws.Range('E29').Value2 = '  -0.40%  '
$ws.Range('E30').Value2 = '  -0.14%  '
$ws.Range('E31').Value2 = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '3.34'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value2 = '  +0.30%  '
$ws.Range('E34').Value2 = '  -1.45%  '
$ws.Range('E35').Value2 = '  -0.36%  '
$ws.Range('E36').Value2 = '  +0.24%  '
$ws.Range('D37').Value2 = '1.131.70'
$ws.Range('E37').Value2 = '  +1.81%  '
$ws.Range('E38').Value2 = '  -2.30%  '
$ws.Range('E39').Value2 = '  -1.48%  '
$ws.Range('E40').Value2 = '  -0.93%  '
$ws.Range('E41').Value2 = '  -0.16%  '
$ws.Range('E42').Value2 = '  -0.33%  '
$ws.Range('E43').Value2 = '  -0.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '100.72'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value2 = '  +1.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.800'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value2 = '  -0.54%  '
$ws.Range('D46').Value2 = '1.769.09'
$ws.Range('E46').Value2 = '  -0.06%  '
$ws.Range('E47').Value2 = '  +0.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '55.43'
$ws.Range('D48').ClearFormats()
$ws.Range('E49').Value2 = '  +0.32%  '
$ws.Range('E50').Value2 = '  -0.52%  '
$ws.Range('E51').Value2 = '  +3.30%  '
